$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, shifting existing rows 21-80 down to 22-81.
$ws.Rows.Item(21).Insert()

# Fill the newly inserted row 21 with the new weekly record.
$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(21, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(21, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(21, 4).Value = 44838
$ws.Cells.Item(21, 5).Value = 15
$ws.Cells.Item(21, 6).Value = 100112012
$ws.Cells.Item(21, 7).Value = "Espinaca"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 250
$ws.Cells.Item(21, 11).Value = 1800
$ws.Cells.Item(21, 12).Value = 2000
$ws.Cells.Item(21, 13).Value = 1900
$ws.Cells.Item(21, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(21, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(21, 16).Value = 633
$ws.Cells.Item(21, 17).Value = 3
$ws.Cells.Item(21, 18).Value = "Hortaliza"
